# Update cached "valor" figures (column E) for the renta_total (rows 176-231)
# and renta_diferencial_precios (rows 292-347) series, years 1963-2018, to
# match the corrected calculation from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comparacion_autores")

$ws.Range("E176").Value = -30.93447024799347
$ws.Range("E177").Value = 70.13063190173146
$ws.Range("E178").Value = 104.9350680740262
$ws.Range("E179").Value = 282.7077624933602
$ws.Range("E180").Value = -61.94102716304131
$ws.Range("E181").Value = -186.2591787303295
$ws.Range("E182").Value = -78.34435110789292
$ws.Range("E183").Value = 91.73485762614818
$ws.Range("E184").Value = 236.7844053683192
$ws.Range("E185").Value = 265.0308929108967
$ws.Range("E186").Value = 643.385941056861
$ws.Range("E187").Value = 1610.870180503204
$ws.Range("E188").Value = 1704.815308530421
$ws.Range("E189").Value = 1466.856164287803
$ws.Range("E190").Value = 1801.472644663213
$ws.Range("E191").Value = 1787.168739195896
$ws.Range("E192").Value = 4947.69068161478
$ws.Range("E193").Value = 6305.362121755735
$ws.Range("E194").Value = 6437.612137187783
$ws.Range("E195").Value = 5761.323671830959
$ws.Range("E196").Value = 4197.004020233644
$ws.Range("E197").Value = 4404.803048276855
$ws.Range("E198").Value = 3606.204688076076
$ws.Range("E199").Value = 2546.964920951899
$ws.Range("E200").Value = 2243.764741000658
$ws.Range("E201").Value = 809.8180589883409
$ws.Range("E202").Value = 3650.642101717535
$ws.Range("E203").Value = 3626.6362723304
$ws.Range("E204").Value = 4742.145264779835
$ws.Range("E205").Value = 3631.224856667284
$ws.Range("E206").Value = 2192.156119528132
$ws.Range("E207").Value = 2038.949164093487
$ws.Range("E208").Value = 2515.947197298002
$ws.Range("E209").Value = 3047.868824296746
$ws.Range("E210").Value = 2753.643386373699
$ws.Range("E211").Value = 1614.181567687775
$ws.Range("E212").Value = 2077.639145860988
$ws.Range("E213").Value = 4433.171726884147
$ws.Range("E214").Value = 3316.748111702572
$ws.Range("E215").Value = 3440.098260076793
$ws.Range("E216").Value = 3044.108096525371
$ws.Range("E217").Value = 4323.037906831347
$ws.Range("E218").Value = 7118.760904441265
$ws.Range("E219").Value = 10357.79276158131
$ws.Range("E220").Value = 11945.32220106251
$ws.Range("E221").Value = 28841.23591783809
$ws.Range("E222").Value = 15603.75819405868
$ws.Range("E223").Value = 21085.2297379748
$ws.Range("E224").Value = 27993.43024244759
$ws.Range("E225").Value = 30146.1357396349
$ws.Range("E226").Value = 24751.17534990787
$ws.Range("E227").Value = 22461.11208471438
$ws.Range("E228").Value = 8216.331359362774
$ws.Range("E229").Value = 3230.964899539834
$ws.Range("E230").Value = 5754.990786191625
$ws.Range("E231").Value = 5489.060682059311
$ws.Range("E292").Value = -30.91298308475931
$ws.Range("E293").Value = 70.10653251581678
$ws.Range("E294").Value = 104.9347648354843
$ws.Range("E295").Value = 282.6894994048039
$ws.Range("E296").Value = -61.94009870151253
$ws.Range("E297").Value = -186.2711967775439
$ws.Range("E298").Value = -78.30309997908293
$ws.Range("E299").Value = 91.73485762614818
$ws.Range("E300").Value = 236.7844053683192
$ws.Range("E301").Value = 265.0308933010269
$ws.Range("E302").Value = 643.3859410568612
$ws.Range("E303").Value = 1609.815740595328
$ws.Range("E304").Value = 1704.401368714782
$ws.Range("E305").Value = 1466.836294605321
$ws.Range("E306").Value = 1801.428227284046
$ws.Range("E307").Value = 1787.092959129895
$ws.Range("E308").Value = 4947.690650757108
$ws.Range("E309").Value = 6305.362121755735
$ws.Range("E310").Value = 6437.612137187783
$ws.Range("E311").Value = 5761.32367183096
$ws.Range("E312").Value = 4197.004020233644
$ws.Range("E313").Value = 4405.041635684676
$ws.Range("E314").Value = 3614.632999643644
$ws.Range("E315").Value = 2548.057723710823
$ws.Range("E316").Value = 2244.040260884682
$ws.Range("E317").Value = 807.5714005251587
$ws.Range("E318").Value = 3664.477951849472
$ws.Range("E319").Value = 3566.328469944585
$ws.Range("E320").Value = 4491.374285651429
$ws.Range("E321").Value = 3274.882673251336
$ws.Range("E322").Value = 1733.911766514953
$ws.Range("E323").Value = 1515.625416887377
$ws.Range("E324").Value = 1647.164991435544
$ws.Range("E325").Value = 1575.717995028259
$ws.Range("E326").Value = 1377.155610781535
$ws.Range("E327").Value = 700.9477028064657
$ws.Range("E328").Value = 1041.087343071962
$ws.Range("E329").Value = 3756.742728648999
$ws.Range("E330").Value = 1648.354237326649
$ws.Range("E331").Value = -462.4130807817629
$ws.Range("E332").Value = 118.7857777366114
$ws.Range("E333").Value = 881.0767653074504
$ws.Range("E334").Value = 3141.631374210432
$ws.Range("E335").Value = 5888.75598252955
$ws.Range("E336").Value = 8549.232893674682
$ws.Range("E337").Value = 23848.25904743331
$ws.Range("E338").Value = 11024.89068590057
$ws.Range("E339").Value = 17196.90897863105
$ws.Range("E340").Value = 24367.32870417406
$ws.Range("E341").Value = 26637.96022537229
$ws.Range("E342").Value = 23122.10508246888
$ws.Range("E343").Value = 20738.2063483997
$ws.Range("E344").Value = 8473.712345993508
$ws.Range("E345").Value = 2254.798464089931
$ws.Range("E346").Value = 4333.686365457092
$ws.Range("E347").Value = 3970.308887884082
